# Update CDA Logical model for ST.r2b
$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" (first sheet) ---
$meta = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3, column B)
$meta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (row 8, column B)
$meta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row before current row 11 ("Description") for the new
# "Jurisdiction" property with an empty value. Copy the formatting from
# the row that will end up right below it so the new row keeps the same
# (border / wrap-text) style used by all the other data rows.
$meta.Rows.Item(11).Insert()
$meta.Range("A12:B12").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""
